$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header row (B1:M1) to new column order ---
$ws.Range("B1").Value = "reddit"
$ws.Range("C1").Value = "patio_lawn_garden"
$ws.Range("D1").Value = "twitter"
$ws.Range("E1").Value = "lidl"
$ws.Range("F1").Value = "automotive"
$ws.Range("G1").Value = "ikea_reviews"
$ws.Range("H1").Value = "luxury_beauty"
$ws.Range("I1").Value = "instant_video"
$ws.Range("J1").Value = "musical_instruments"
$ws.Range("K1").Value = "office_products"
$ws.Range("L1").Value = "hotel"
$ws.Range("M1").Value = "drugs"

# --- Move existing SVM row (row 7) to row 8, preserving its style, before overwriting rows ---
$ws.Range("A7").Copy($ws.Range("A8"))

# --- Write row labels and data rows 2-8 in new order/values ---
$ws.Range("A2").Value = "ComplementNB"
$ws.Range("B2").Value = 66.95999999999999
$ws.Range("C2").Value = 50.43
$ws.Range("D2").Value = 78.01000000000001
$ws.Range("E2").Value = 64.09999999999999
$ws.Range("F2").Value = 60.71
$ws.Range("G2").Value = 66.39
$ws.Range("H2").Value = 67.42
$ws.Range("I2").Value = 56.23
$ws.Range("J2").Value = 58.72
$ws.Range("K2").Value = 56.27
$ws.Range("L2").Value = 50.16
$ws.Range("M2").Value = 37.83

$ws.Range("A3").Value = "Decision Tree"
$ws.Range("B3").Value = 60.37
$ws.Range("C3").Value = 42.8
$ws.Range("D3").Value = 67.69
$ws.Range("E3").Value = 58.7
$ws.Range("F3").Value = 55.36
$ws.Range("G3").Value = 60.76
$ws.Range("H3").Value = 64.52
$ws.Range("I3").Value = 45.84
$ws.Range("J3").Value = 56.07
$ws.Range("K3").Value = 46.38
$ws.Range("L3").Value = 39.48
$ws.Range("M3").Value = 33.99

$ws.Range("A4").Value = "LR"
$ws.Range("B4").Value = 70.39
$ws.Range("C4").Value = 54.35
$ws.Range("D4").Value = 80.56
$ws.Range("E4").Value = 66.41
$ws.Range("F4").Value = 62.39
$ws.Range("G4").Value = 71.14
$ws.Range("H4").Value = 73.23
$ws.Range("I4").Value = 59.57
$ws.Range("J4").Value = 62.56
$ws.Range("K4").Value = 61.17
$ws.Range("L4").Value = 52.19
$ws.Range("M4").Value = 42.32

$ws.Range("A5").Value = "MultinomialNB"
$ws.Range("B5").Value = 67.16
$ws.Range("C5").Value = 47.1
$ws.Range("D5").Value = 73.43000000000001
$ws.Range("E5").Value = 64.54000000000001
$ws.Range("F5").Value = 56.83
$ws.Range("G5").Value = 69.14
$ws.Range("H5").Value = 63.67
$ws.Range("I5").Value = 51.3
$ws.Range("J5").Value = 56.52
$ws.Range("K5").Value = 50.13
$ws.Range("L5").Value = 49.11
$ws.Range("M5").Value = 35.52

$ws.Range("A6").Value = "RF"
$ws.Range("B6").Value = 61.4
$ws.Range("C6").Value = 37.53
$ws.Range("D6").Value = 56.07
$ws.Range("E6").Value = 51.33
$ws.Range("F6").Value = 54.25
$ws.Range("G6").Value = 60.46
$ws.Range("H6").Value = 52.91
$ws.Range("I6").Value = 41.86
$ws.Range("J6").Value = 55.15
$ws.Range("K6").Value = 41.09
$ws.Range("L6").Value = 31.84
$ws.Range("M6").Value = 17.92

$ws.Range("A7").Value = "setfit"
$ws.Range("B7").Value = 58.99
$ws.Range("C7").Value = 63.57
$ws.Range("D7").Value = 73.31999999999999
$ws.Range("E7").Value = 53.21
$ws.Range("F7").Value = 59.24
$ws.Range("G7").Value = 67.37
$ws.Range("H7").Value = 69.56
$ws.Range("I7").Value = 67.06
$ws.Range("J7").Value = 73.31999999999999
$ws.Range("K7").Value = 68.11
$ws.Range("L7").Value = 71.06999999999999
$ws.Range("M7").Value = 58.61

$ws.Range("A8").Value = "SVM"
$ws.Range("B8").ClearContents()
$ws.Range("C8").ClearContents()
$ws.Range("D8").ClearContents()
$ws.Range("E8").Value = 67.38
$ws.Range("F8").Value = 64.64
$ws.Range("G8").Value = 70.33
$ws.Range("H8").Value = 75.09999999999999
$ws.Range("I8").Value = 59.99
$ws.Range("J8").Value = 61.56
$ws.Range("K8").ClearContents()
$ws.Range("L8").Value = 53.67
$ws.Range("M8").Value = 46.73
